# New weekly price-report row for "Cebollín baby" (Agrícola del Norte S.A. de Arica).
# A new record is inserted at row 71, pushing the existing rows 71-88 down to 72-89.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 71 - this shifts rows 71:88 down to 72:89
# and keeps all other data (and row formatting) intact.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with this week's data.
$ws.Range("A71").Value = 1
$ws.Range("B71").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C71").Value = "Arica y Parinacota"
$ws.Range("D71").Value = 44736
$ws.Range("E71").Value = 15
$ws.Range("F71").Value = 100112038
$ws.Range("G71").Value = "Cebollín baby"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 270
$ws.Range("K71").Value = 2800
$ws.Range("L71").Value = 3000
$ws.Range("M71").Value = 2900
$ws.Range("N71").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 1450
$ws.Range("Q71").Value = 2
$ws.Range("R71").Value = "Hortaliza"
